$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63..105 down to 64..106
$ws.Rows.Item(63).Insert()

# Copy the (now shifted-down) old row 63 content, which now lives at row 64,
# into the new blank row 63, then adjust the few cells that actually changed.
$srcRow = 64
$dstRow = 63
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($dstRow, $c).Value = $ws.Cells.Item($srcRow, $c).Value2
}

# Match the date cell's number format/style to the rest of the date column (style index 2 / numFmtId 165)
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat

# Now set the new values for the inserted row (D, M, N, O, P, S)
$ws.Cells.Item($dstRow, 4).Value = 44596   # Fecha
$ws.Cells.Item($dstRow, 13).Value = 150    # Volumen
$ws.Cells.Item($dstRow, 14).Value = 7000   # Precio minimo
$ws.Cells.Item($dstRow, 15).Value = 7000   # Precio maximo
$ws.Cells.Item($dstRow, 16).Value = 7000   # Precio promedio ponderado
$ws.Cells.Item($dstRow, 19).Value = 1750   # Precio $/Kg
